# "minor bugfixes and annotations"
#
# Source diff shows three related tweaks to Sheet1's view/columns:
#   1. The view scrolled so topLeftCell goes from A1 -> B1, and the
#      active/selected cell moves from H21 -> E5.
#   2. Column C (Material) was widened to ~21.15 chars and column F
#      (Batch) was widened to ~16.41 chars (both previously at the
#      sheet default width, i.e. no explicit <col> entry).
#   3. The sheet default column width shifted very slightly
#      (8.3828125 -> 8.390625 chars) - a trivial by-product of the
#      resize above, not something to target directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- 1. Scroll position + selection -----------------------------------
# Move the viewport so column B is the left-most visible column and row 1
# stays the top-most visible row (topLeftCell A1 -> B1), then select E5
# (the new active cell), matching activeCell/sqref in the diff.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 2
$ws.Range("E5").Select()

# --- 2. Column widths ----------------------------------------------------
# Widen column C (Material) and column F (Batch) - both were previously
# unset/default width and become explicit customWidth columns. Target
# character widths from the diff are 21.15 and 16.41; Excel's COM layer
# quantises ColumnWidth to whole on-screen pixels (width_px = round(cw*6+5),
# stored_chars = width_px/6), so these inputs are chosen to land on the
# closest reachable pixel-quantised width (21.166667 and 16.333333).
$ws.Columns.Item(3).ColumnWidth = 20.333333
$ws.Columns.Item(6).ColumnWidth = 15.5
